# Update "Local" and "Donor" sheets to 2017 actual figures.
# Rows 2-82 get new district order + values; rows 83-112 are removed
# (the dataset now covers 81 districts instead of 111).

$wb = $excel.ActiveWorkbook

$localData = @(
    @("Kalangala", 4.541639438),
    @("Kiboga", 5.203772724),
    @("Luwero", 0.916304667),
    @("Masaka", 14.18571894),
    @("Mpigi", 4.139156418),
    @("Mubende", 3.538695403),
    @("Mukono", 4.088950403),
    @("Rakai", 4.085117158),
    @("Kayunga", 4.052338098),
    @("Wakiso", 16.87423146),
    @("Lyantonde", 5.983323335),
    @("Mityana", 4.756896631),
    @("Buikwe", 28.25647561),
    @("Butambala", 0.732599047),
    @("Buvuma", 4.0636447),
    @("Gomba", 3.901773227),
    @("Kalungu", 3.486128824),
    @("Lwengo", 3.861227727),
    @("Busia", 0.978339196),
    @("Jinja", 38.83940352),
    @("Kamuli", 2.787675764),
    @("Kapchorwa", 2.72585879),
    @("Kumi", 12.91936959),
    @("Mbale", 1.692963805),
    @("Soroti", 6.241647428),
    @("Tororo", 6.474493141),
    @("Mayuge", 1.946048319),
    @("Sironko", 1.813154734),
    @("Amuria", 2.560169377),
    @("Budaka", 1.298996415),
    @("Bududa", 1.798340469),
    @("Butaleja", 0.978339196),
    @("Kaliro", 2.329644826),
    @("Namutumba", 2.297467716),
    @("Buyende", 1.674426788),
    @("Kibuku", 0.46408253),
    @("Kween", 1.785610351),
    @("Luuka", 0.489512484),
    @("Ngora", 4.441923821),
    @("Serere", 4.405825572),
    @("Adjumani", 1.134584252),
    @("Arua", 4.276503322),
    @("Gulu", 2.264433981),
    @("Kotido", 2.242613489),
    @("Lira", 4.014831409),
    @("Moyo", 2.998551509),
    @("Nebbi", 9.326400989),
    @("Nakapiripirit", 1.400130352),
    @("Yumbe", 1.806781911),
    @("Amolatar", 2.198432456),
    @("Amuru", 2.424495698),
    @("Dokolo", 0.817709415),
    @("Koboko", 4.184856282),
    @("Maracha", 1.216133273),
    @("Oyam", 2.063990747),
    @("Alebtong", 2.069804033),
    @("Amudat", 2.069099584),
    @("Kole", 1.74498138),
    @("Zombo", 7.419495079),
    @("Bushenyi", 12.34923114),
    @("Hoima", 6.558630433),
    @("Kabale", 8.995898879),
    @("Kabarole", 8.373917445),
    @("Kasese", 12.83532329),
    @("Kibaale", 3.103542477),
    @("Kisoro", 27.85784585),
    @("Masindi", 16.21956159),
    @("Mbarara", 13.83125798),
    @("Ntungamo", 4.439137549),
    @("Rukungiri", 1.799115782),
    @("Kamwenge", 2.611909214),
    @("Kanungu", 2.078787564),
    @("Kyenjojo", 3.556581677),
    @("Ibanda", 7.25330341),
    @("Isingiro", 3.898251801),
    @("Kiruhura", 6.671384551),
    @("Kiryandongo", 5.594435146),
    @("Kyegegwa", 1.547664591),
    @("Mitooma", 1.728132866),
    @("Rubirizi", 2.853950334),
    @("Sheema", 1.942748259)
)

$donorData = @(
    @("Kalangala", 53.31404308),
    @("Kiboga", 1.408770599),
    @("Luwero", 1.774007196),
    @("Masaka", $null),
    @("Mpigi", 10.09650351),
    @("Mubende", 3.610177056),
    @("Mukono", 2.841294712),
    @("Rakai", 3.113850986),
    @("Kayunga", 3.627060691),
    @("Wakiso", 1.880065746),
    @("Lyantonde", 2.87234637),
    @("Mityana", 0.533288975),
    @("Buikwe", $null),
    @("Butambala", $null),
    @("Buvuma", 6.696821022),
    @("Gomba", 2.176079603),
    @("Kalungu", 3.256363829),
    @("Lwengo", 2.27902315),
    @("Busia", 1.625933957),
    @("Jinja", $null),
    @("Kamuli", 1.43976639),
    @("Kapchorwa", 0.496146093),
    @("Kumi", $null),
    @("Mbale", 5.820147339),
    @("Soroti", 3.243644532),
    @("Tororo", 1.438782402),
    @("Mayuge", 14.38747837),
    @("Sironko", 1.89660106),
    @("Amuria", 1.86650182),
    @("Budaka", 0.546306064),
    @("Bududa", 3.125688367),
    @("Butaleja", 1.625933957),
    @("Kaliro", 2.675016686),
    @("Namutumba", 1.63010458),
    @("Buyende", 1.163413673),
    @("Kibuku", $null),
    @("Kween", 1.182252391),
    @("Luuka", 6.527998504),
    @("Ngora", $null),
    @("Serere", 0.823012001),
    @("Adjumani", 14.97864358),
    @("Arua", 1.555069352),
    @("Gulu", 1.789579143),
    @("Kotido", 17.88763843),
    @("Lira", 1.580202818),
    @("Moyo", 23.81550264),
    @("Nebbi", $null),
    @("Nakapiripirit", 16.35479354),
    @("Yumbe", 8.789180489),
    @("Amolatar", 49.44722943),
    @("Amuru", 2.355256429),
    @("Dokolo", 1.812884594),
    @("Koboko", 13.72285981),
    @("Maracha", 2.320414688),
    @("Oyam", 12.9731581),
    @("Alebtong", 0.957801758),
    @("Amudat", 8.242166581),
    @("Kole", 0.596423192),
    @("Zombo", 3.011438326),
    @("Bushenyi", $null),
    @("Hoima", 3.440094174),
    @("Kabale", $null),
    @("Kabarole", 1.106719661),
    @("Kasese", 1.883875474),
    @("Kibaale", 4.388483053),
    @("Kisoro", $null),
    @("Masindi", $null),
    @("Mbarara", $null),
    @("Ntungamo", 10.14051046),
    @("Rukungiri", 0.882331804),
    @("Kamwenge", 1.056106401),
    @("Kanungu", 4.949204454),
    @("Kyenjojo", 2.894125184),
    @("Ibanda", 4.07684759),
    @("Isingiro", 2.644795613),
    @("Kiruhura", 5.558078944),
    @("Kiryandongo", 0.794105668),
    @("Kyegegwa", 5.124140376),
    @("Mitooma", 0.230206164),
    @("Rubirizi", 3.415749762),
    @("Sheema", 1.157516618)
)

$wsLocal = $wb.Worksheets.Item("Local")
$r = 2
foreach ($row in $localData) {
    $wsLocal.Cells.Item($r, 1).Value = $row[0]
    $wsLocal.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
$wsLocal.Range("A83:B112").EntireRow.Delete()

$wsDonor = $wb.Worksheets.Item("Donor")
$r = 2
foreach ($row in $donorData) {
    $wsDonor.Cells.Item($r, 1).Value = $row[0]
    $wsDonor.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
$wsDonor.Range("A83:B112").EntireRow.Delete()
